$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing data (rows 1-162) shifts down to rows 2-163.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Município"
$ws.Range("B1").Value = "Casos"
$ws.Range("C1").Value = "Óbitos"

# Build the header formatting (bold, thin box border, centered/top-aligned) on an unused
# helper cell, then copy/paste the formatting onto the header row in a single operation so
# the style table only gains one combined cell format (matches a single new cellXf).
$helper = $ws.Range("Z1")
$helper.Font.Bold = $true
$helper.Borders.LineStyle = 1
$helper.HorizontalAlignment = -4108
$helper.VerticalAlignment = -4160

$helper.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$helper.Clear()

# Append two new trailing rows.
$ws.Range("A164").Value = "outros estados"
$ws.Range("B164").Value = 43

$ws.Range("A165").Value = "outros paises"
$ws.Range("B165").Value = 39
